# The document was edited in MS Word: the grammar checker flagged every
# '({ width' fragment (added w:proofErr gramStart/gramEnd bookmarks around it)
# and the spell checker flagged 'sample.svg' in the last placeholder
# (w:proofErr spellStart/spellEnd). The actual content edit is adding SVG
# thumbnail support: './sample.svg' => './sample.svg', thumbnail: './sample.png'
#
# We reproduce this by rebuilding each affected paragraph's WordprocessingML
# (preserving original run-level w:rsidR markers) and injecting it with
# Range.InsertXML, which replaces exactly the target paragraph's contents.
$d = $word.ActiveDocument

# Paragraph 2: verify it is the expected placeholder, then replace
# its markup with the proofed / edited version.
$p2 = $d.Paragraphs[2]
$r2 = $p2.Range
if ($r2.Text -notlike "*+++IMAGE ({ width: 3, height: 3, path: './sample.png' })+++*") {
    throw "Paragraph 2 did not contain expected text: $($r2.Text)"
}
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="425B0F31" w14:textId="345E38FE" w:rsidR="004E34BD" w:rsidRDefault="00904BD3"><w:r><w:t xml:space="preserve">+++IMAGE </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00295C1E"><w:t>(</w:t></w:r><w:r><w:t>{ width</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>: 3, height: 3, path: '</w:t></w:r><w:r w:rsidR="00DF48F9"><w:t>./</w:t></w:r><w:r w:rsidR="004C35D4"><w:t>sample</w:t></w:r><w:r w:rsidR="00DF48F9"><w:t>.png</w:t></w:r><w:r><w:t>' }</w:t></w:r><w:r w:rsidR="00295C1E"><w:t>)</w:t></w:r><w:r><w:t>+++</w:t></w:r></w:p>
'@
[void]$r2.InsertXML($xml2)

# Paragraph 3: verify it is the expected placeholder, then replace
# its markup with the proofed / edited version.
$p3 = $d.Paragraphs[3]
$r3 = $p3.Range
if ($r3.Text -notlike "*+++IMAGE ({ width: 3, height: 3, path: './sample.jpg' })+++*") {
    throw "Paragraph 3 did not contain expected text: $($r3.Text)"
}
$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7F1455EC" w14:textId="4F8D11D8" w:rsidR="004C35D4" w:rsidRDefault="004C35D4" w:rsidP="004C35D4"><w:r><w:t xml:space="preserve">+++IMAGE </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>({ width</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>: 3, height: 3, path: './sample.jpg</w:t></w:r><w:r w:rsidR="00292508"><w:t>'</w:t></w:r><w:r w:rsidR="00634825"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>})+++</w:t></w:r></w:p>
'@
[void]$r3.InsertXML($xml3)

# Paragraph 4: verify it is the expected placeholder, then replace
# its markup with the proofed / edited version.
$p4 = $d.Paragraphs[4]
$r4 = $p4.Range
if ($r4.Text -notlike "*+++IMAGE ({ width: 3, height: 3, path: './sample.jpeg' })+++*") {
    throw "Paragraph 4 did not contain expected text: $($r4.Text)"
}
$xml4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2BF4E645" w14:textId="6D32B310" w:rsidR="008A4E56" w:rsidRDefault="008A4E56" w:rsidP="008A4E56"><w:r><w:t xml:space="preserve">+++IMAGE </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>({ width</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>: 3, height: 3, path: './sample.jpeg' })+++</w:t></w:r></w:p>
'@
[void]$r4.InsertXML($xml4)

# Paragraph 5: verify it is the expected placeholder, then replace
# its markup with the proofed / edited version.
$p5 = $d.Paragraphs[5]
$r5 = $p5.Range
if ($r5.Text -notlike "*+++IMAGE ({ width: 3, height: 3, path: './sample.gif' })+++*") {
    throw "Paragraph 5 did not contain expected text: $($r5.Text)"
}
$xml5 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4A31339F" w14:textId="6FFC37EA" w:rsidR="004C35D4" w:rsidRDefault="004C35D4" w:rsidP="004C35D4"><w:r><w:t xml:space="preserve">+++IMAGE </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>({ width</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>: 3, height: 3, path: './sample.gif' })+++</w:t></w:r></w:p>
'@
[void]$r5.InsertXML($xml5)

# Paragraph 6: verify it is the expected placeholder, then replace
# its markup with the proofed / edited version.
$p6 = $d.Paragraphs[6]
$r6 = $p6.Range
if ($r6.Text -notlike "*+++IMAGE ({ width: 3, height: 3, path: './sample.bmp' })+++*") {
    throw "Paragraph 6 did not contain expected text: $($r6.Text)"
}
$xml6 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5DE41E8A" w14:textId="78ADECF5" w:rsidR="000763F6" w:rsidRDefault="000763F6" w:rsidP="000763F6"><w:r><w:t xml:space="preserve">+++IMAGE </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>({ width</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>: 3, height: 3, path: './sample.bmp' })+++</w:t></w:r></w:p>
'@
[void]$r6.InsertXML($xml6)

# Paragraph 7: verify it is the expected placeholder, then replace
# its markup with the proofed / edited version.
$p7 = $d.Paragraphs[7]
$r7 = $p7.Range
if ($r7.Text -notlike "*+++IMAGE ({ width: 3, height: 3, path: './sample.svg' })+++*") {
    throw "Paragraph 7 did not contain expected text: $($r7.Text)"
}
$xml7 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="057C6856" w14:textId="6690D083" w:rsidR="008A3722" w:rsidRDefault="008A3722" w:rsidP="008A3722"><w:r><w:t xml:space="preserve">+++IMAGE </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>({ width</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>: 3, height: 3, path: './</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sample.</w:t></w:r><w:r w:rsidR="00B601A3"><w:t>svg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>'</w:t></w:r><w:r><w:t>,</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> thumbnail:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">'./sample.png' </w:t></w:r><w:r><w:t>})+++</w:t></w:r></w:p>
'@
[void]$r7.InsertXML($xml7)

